$wb = $excel.ActiveWorkbook

# --- Insert the new "WISdom results" sheet right before "BTaDLP" (i.e. after "Texas Notes") ---
$texasNotes = $wb.Worksheets.Item("Texas Notes")
$wisdom = $wb.Worksheets.Add($null, $texasNotes)
$wisdom.Name = "WISdom results"

# Content: a labeled average transmission-loss figure used by BTaDLP.
$wisdom.Range("A1").Value = "Average transmission losses across all TX regions"
$wisdom.Range("A1").WrapText = $true
$wisdom.Range("B1").Value = 0.06
$wisdom.Rows.Item(1).RowHeight = 28.5
$wisdom.Columns.Item(1).ColumnWidth = 24.6
$wisdom.Range("B2").Select() | Out-Null

# --- Point BTaDLP's T&D loss row at the new WISdom results figure instead of the old TREND()/Texas Notes formulas ---
$btadlp = $wb.Worksheets.Item("BTaDLP")
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z", `
          "AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK")
foreach ($col in $cols) {
    $btadlp.Range("$col`2").Formula = "='WISdom results'!`$B`$1"
}

# Reflect the new selection left behind on the BTaDLP sheet.
$btadlp.Range("B2:AK2").Select() | Out-Null

# Leave the workbook with "About" as the active tab (sheet 1).
$about = $wb.Worksheets.Item("About")
$about.Activate()
